$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-16 02:48:44'
$ws.Range('E3').Value = '2026-02-16 02:48:46'
$ws.Range('I3').Value = '0.4 mm'
$ws.Range('E4').Value = '2026-02-16 02:48:48'
$ws.Range('H4').Value = "'70%"
$ws.Range('G4').Copy()
$ws.Range('H4').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('J4').Value = '1014.8 hPa'
$ws.Range('O4').Value = '10.4 °C'
$ws.Range('E5').Value = '2026-02-16 02:48:51'
$ws.Range('I5').Value = '1.3 mm'
$ws.Range('E6').Value = '2026-02-16 02:48:53'
$ws.Range('N6').Value = '6.3 °C 2:05 TU'
$ws.Range('E7').Value = '2026-02-16 02:48:56'
$ws.Range('J7').Value = '1015.3 hPa'
$ws.Range('E8').Value = '2026-02-16 02:48:58'
$ws.Range('H8').Value = "'69%"
$ws.Range('G8').Copy()
$ws.Range('H8').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('J8').Value = '1015.2 hPa'
$ws.Range('L8').Value = '56.5 km/h - 300º 2:22 TU'
$ws.Range('E9').Value = '2026-02-16 02:49:01'
$ws.Range('H9').Value = "'93%"
$ws.Range('G9').Copy()
$ws.Range('H9').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('O9').Value = '5.7 °C'
$ws.Range('E10').Value = '2026-02-16 02:49:03'
$ws.Range('L10').Value = '5.0 km/h - 92º 2:26 TU'
$ws.Range('E11').Value = '2026-02-16 02:49:05'
$ws.Range('E12').Value = '2026-02-16 02:49:08'
$ws.Range('N12').Value = '4.8 °C 2:02 TU'
$ws.Range('O12').Value = '5.7 °C'
$ws.Range('E13').Value = '2026-02-16 02:49:10'
$ws.Range('J13').Value = '1018.6 hPa'
$ws.Range('L13').Value = '11.5 km/h - 353º 2:19 TU'
$ws.Range('O13').Value = '1.6 °C'
$ws.Range('E14').Value = '2026-02-16 02:49:11'
$ws.Range('M14').Value = '13.1 °C 2:25 TU'
$ws.Range('O14').Value = '12.5 °C'
$ws.Range('E15').Value = '2026-02-16 02:49:13'
$ws.Range('H15').Value = "'88%"
$ws.Range('G15').Copy()
$ws.Range('H15').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('E16').Value = '2026-02-16 02:49:14'
$ws.Range('M16').Value = '-0.6 °C 2:29 TU'
$ws.Range('E17').Value = '2026-02-16 02:49:15'
$ws.Range('H17').Value = "'60%"
$ws.Range('G17').Copy()
$ws.Range('H17').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('L17').Value = '35.3 km/h - 258º 2:12 TU'
$ws.Range('N17').Value = '5.0 °C 2:26 TU'
$ws.Range('O17').Value = '5.5 °C'
$ws.Range('E18').Value = '2026-02-16 02:49:16'
$ws.Range('H18').Value = "'97%"
$ws.Range('G18').Copy()
$ws.Range('H18').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('N18').Value = '3.2 °C 2:09 TU'
$ws.Range('O18').Value = '4.6 °C'
$ws.Range('E19').Value = '2026-02-16 02:49:17'
$ws.Range('N19').Value = '2.9 °C 2:19 TU'
$ws.Range('O19').Value = '3.5 °C'
$ws.Range('E20').Value = '2026-02-16 02:49:18'
$ws.Range('N20').Value = '-1.2 °C 2:11 TU'
$ws.Range('O20').Value = '-0.9 °C'
$ws.Range('E21').Value = '2026-02-16 02:49:19'
$ws.Range('N21').Value = '4.6 °C 2:17 TU'
$ws.Range('O21').Value = '5.3 °C'
$ws.Range('E22').Value = '2026-02-16 02:49:20'
$ws.Range('N22').Value = '-6.5 °C 2:29 TU'
$ws.Range('E23').Value = '2026-02-16 02:49:21'
$ws.Range('I23').Value = '0.5 mm'
$ws.Range('E24').Value = '2026-02-16 02:49:22'
$ws.Range('H24').Value = "'69%"
$ws.Range('G24').Copy()
$ws.Range('H24').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('J24').Value = '1018.4 hPa'
$ws.Range('E25').Value = '2026-02-16 02:49:25'
$ws.Range('H25').Value = "'71%"
$ws.Range('G25').Copy()
$ws.Range('H25').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('O25').Value = '0.8 °C'
$ws.Range('E26').Value = '2026-02-16 02:49:27'
$ws.Range('E27').Value = '2026-02-16 02:49:29'
$ws.Range('H27').Value = "'74%"
$ws.Range('G27').Copy()
$ws.Range('H27').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('E28').Value = '2026-02-16 02:49:32'
$ws.Range('J28').Value = '1016.2 hPa'
$ws.Range('N28').Value = '2.6 °C 2:06 TU'
$ws.Range('O28').Value = '3.3 °C'
$ws.Range('E29').Value = '2026-02-16 02:49:34'
$ws.Range('E30').Value = '2026-02-16 02:49:37'
$ws.Range('J30').Value = '1015.0 hPa'
$ws.Range('M30').Value = '7.3 °C 2:27 TU'
$ws.Range('N30').Value = '6.4 °C 2:05 TU'
$ws.Range('E31').Value = '2026-02-16 02:49:39'
$ws.Range('H31').Value = "'56%"
$ws.Range('G31').Copy()
$ws.Range('H31').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('J31').Value = '1013.5 hPa'
$ws.Range('L31').Value = '74.5 km/h - 349º 2:20 TU'
$ws.Range('O31').Value = '14.0 °C'
$ws.Range('E32').Value = '2026-02-16 02:49:42'
$ws.Range('O32').Value = '6.2 °C'
$ws.Range('E33').Value = '2026-02-16 02:49:44'
$ws.Range('H33').Value = "'68%"
$ws.Range('G33').Copy()
$ws.Range('H33').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('J33').Value = '1015.9 hPa'
$ws.Range('N33').Value = '4.4 °C 2:29 TU'
$ws.Range('O33').Value = '5.5 °C'
$ws.Range('E34').Value = '2026-02-16 02:49:47'
$ws.Range('H34').Value = "'64%"
$ws.Range('G34').Copy()
$ws.Range('H34').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('N34').Value = '2.5 °C 2:03 TU'
$ws.Range('O34').Value = '3.4 °C'
$ws.Range('E35').Value = '2026-02-16 02:49:49'
$ws.Range('M35').Value = '6.9 °C 2:29 TU'
$ws.Range('E36').Value = '2026-02-16 02:49:52'
$ws.Range('H36').Value = "'85%"
$ws.Range('G36').Copy()
$ws.Range('H36').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('J36').Value = '1014.8 hPa'
$ws.Range('N36').Value = '6.3 °C 2:28 TU'
$ws.Range('O36').Value = '7.4 °C'
$ws.Range('E37').Value = '2026-02-16 02:49:54'
$ws.Range('J37').Value = '1018.3 hPa'
$ws.Range('N37').Value = '1.8 °C 2:09 TU'
$ws.Range('O37').Value = '2.2 °C'
$ws.Range('E38').Value = '2026-02-16 02:49:57'
$ws.Range('H38').Value = "'91%"
$ws.Range('G38').Copy()
$ws.Range('H38').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('O38').Value = '5.8 °C'
$ws.Range('E39').Value = '2026-02-16 02:49:59'
$ws.Range('H39').Value = "'75%"
$ws.Range('G39').Copy()
$ws.Range('H39').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('O39').Value = '-0.2 °C'
$ws.Range('E40').Value = '2026-02-16 02:50:02'
$ws.Range('N40').Value = '3.2 °C 2:29 TU'
$ws.Range('O40').Value = '3.5 °C'
$ws.Range('E41').Value = '2026-02-16 02:50:04'
$ws.Range('N41').Value = '14.4 °C 2:25 TU'
$ws.Range('E42').Value = '2026-02-16 02:50:07'
$ws.Range('O42').Value = '6.4 °C'
$ws.Range('E43').Value = '2026-02-16 02:50:09'
$ws.Range('N43').Value = '2.6 °C 2:25 TU'
$ws.Range('O43').Value = '3.8 °C'
$ws.Range('E44').Value = '2026-02-16 02:50:11'
$ws.Range('H44').Value = "'88%"
$ws.Range('G44').Copy()
$ws.Range('H44').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('I44').Value = '0.9 mm'
$ws.Range('M44').Value = '0.1 °C 2:04 TU'
$ws.Range('O44').Value = '-0.5 °C'
$ws.Range('E45').Value = '2026-02-16 02:50:14'
$ws.Range('I45').Value = '0.9 mm'
$ws.Range('J45').Value = '1020.3 hPa'
$ws.Range('L45').Value = '4.7 km/h - 170º 2:06 TU'
$ws.Range('M45').Value = '3.4 °C 2:29 TU'
$ws.Range('N45').Value = '3.1 °C 2:20 TU'
$ws.Range('E46').Value = '2026-02-16 02:50:16'
$ws.Range('H46').Value = "'61%"
$ws.Range('G46').Copy()
$ws.Range('H46').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('J46').Value = '1018.9 hPa'
$ws.Range('M46').Value = '13.0 °C 2:10 TU'
$ws.Range('O46').Value = '12.2 °C'
